$wb = $excel.ActiveWorkbook

# --- Sheet1: Statistics ---
$ws1 = $wb.Worksheets.Item("Statistics")

$ws1.Range("A2").Value = "2024-08-25 21:16:20"
$ws1.Range("B2").Value = 42.4964350850333
$ws1.Range("C2").Value = 5

$ws1.Range("A3").Value = "2024-08-25 21:16:22"
$ws1.Range("B3").Value = 42.33126581995717
$ws1.Range("C3").Value = 7

$ws1.Range("A4").Value = "2024-08-25 21:16:24"
$ws1.Range("B4").Value = 43.87835996264035
$ws1.Range("C4").Value = 11

$ws1.Range("A5").Value = "2024-08-25 21:16:26"
$ws1.Range("B5").Value = 43.00827997015465
$ws1.Range("C5").Value = 11

$ws1.Range("A6").Value = "2024-08-25 21:16:28"
$ws1.Range("B6").Value = 41.63382863258526
$ws1.Range("C6").Value = 15

$ws1.Range("A7").Value = "2024-08-25 21:16:30"
$ws1.Range("B7").Value = 31.75084319358936
$ws1.Range("C7").Value = 14

# --- Sheet2: Accidents ---
$ws2 = $wb.Worksheets.Item("Accidents")

$ws2.Range("A2").Value = "2024-08-25 21:16:30"
$ws2.Range("B2").Value = "Car and Car"
$ws2.Range("C2").Value = "35.89 and 39.31"
$ws2.Range("D2").Value = 1
